# Updates the crypto price/volume snapshot on Sheet1 (columns D = Price, E = Volume(1h))
# to the latest scraped values, cell by cell, matching the upstream GitHub Actions commit.
#
# Column D cells that are themselves unambiguous decimal numbers (e.g. "573.83") would be
# auto-coerced to the Number type by Excel's normal text-to-value parsing, which would both
# change the cell's stored type away from Text and silently drop significant trailing zeros
# (e.g. "10.00" -> 10). Those specific cells are written with a leading apostrophe, Excel's
# standard "force text" quote-prefix, so they stay text cells with the exact literal string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '62.742.47'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '2.464.20'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'573.83"
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = "'147.63"
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').Value = "'5.30"
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').Value = "'29.10"
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').Value = "'0.0000177"
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '62.663.60'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').Value = '2.454.58'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = "'7.92"
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('D19').Value = "'10.91"
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('D20').Value = "'326.21"
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = "'10.00"
$ws.Range('E24').Value = '  +11.53%  '
$ws.Range('D25').Value = "'65.40"
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('D26').Value = "'643.16"
$ws.Range('E26').Value = '  -3.11%  '
$ws.Range('D28').Value = '0.0₃0974'
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('E29').Value = '  -15.41%  '
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').Value = "'7.95"
$ws.Range('E31').Value = '  -2.86%  '
$ws.Range('E32').Value = '  -3.47%  '
$ws.Range('E33').Value = '  -4.06%  '
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('D36').Value = "'4.75"
$ws.Range('E36').Value = '  -0.87%  '
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('D38').Value = "'150.75"
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('D40').Value = "'5.31"
$ws.Range('E40').Value = '  -3.48%  '
$ws.Range('D41').Value = "'2.73"
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('D43').Value = '0.0₆0313'
$ws.Range('E43').Value = '  -10.39%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = "'153.17"
$ws.Range('E45').Value = '  +4.58%  '
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('E47').Value = '  -1.23%  '
$ws.Range('D48').Value = "'20.37"
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('D49').Value = "'0.607"
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = "'0.0508"
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('E51').Value = '  -1.36%  '
